# Auto commit at 2025-09-29 15:51:07.95
#
# Updates the "Metrics" sheet figures (B2:B13) with refreshed values. The
# "today" sheet pulls these same figures via formulas (=Metrics!Bn) so its
# cached results are recalculated automatically. Also refreshes the sheet
# selection / active-tab state left behind by the editor.

$wb = $excel.ActiveWorkbook

# --- Update the Metrics figures -------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 417310.57
$metrics.Range("B3").Value  = 337330.87
$metrics.Range("B4").Value  = 131646.82999999999
$metrics.Range("B5").Value  = 16586
$metrics.Range("B6").Value  = 4336561.49
$metrics.Range("B7").Value  = 3664858.35
$metrics.Range("B8").Value  = 1261012.51
$metrics.Range("B9").Value  = 167746
$metrics.Range("B10").Value = 32801885.280000001
$metrics.Range("B11").Value = 30940079.91
$metrics.Range("B12").Value = 11542721.42
$metrics.Range("B13").Value = 1265373

# --- Refresh the view/selection state left in the file ---------------------------
# The "today" sheet was scrolled and had B15 selected before the user moved on.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("B15").Select()

# Finally the Metrics sheet became the active tab, with D14 selected.
$metrics.Activate()
$metrics.Range("D14").Select()
